$ErrorActionPreference = "Stop"

$p = $ppt.ActivePresentation

# Locate the paragraph containing the sentence that needs to be split into
# multiple runs: "Certain towers have different affects where they are"
# (the target edit also fixes "affects" -> "effects").
$needle = "Certain towers have different affects where they are"

$targetPara = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text.Contains($needle)) {
                $paraCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $paraCount; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    if ($para.Text.Contains($needle)) {
                        $targetPara = $para
                    }
                }
            }
        }
    }
}

if ($targetPara -eq $null) {
    throw "Could not find target paragraph"
}

# Step 1: split out "different " into its own run (re-typing the same word
# is what produces a dedicated run boundary without touching the rest of
# the sentence).
$text1 = $targetPara.Text
$pos1 = $text1.IndexOf("different")
if ($pos1 -lt 0) {
    throw "Could not find 'different' in target paragraph"
}
$start1 = $pos1 + 1
$range1 = $targetPara.Characters($start1, 10)
$range1.Text = "different "

# Step 2: split out "affects " into its own run and correct it to
# "effects ".
$text2 = $targetPara.Text
$pos2 = $text2.IndexOf("affects")
if ($pos2 -lt 0) {
    throw "Could not find 'affects' in target paragraph"
}
$start2 = $pos2 + 1
$range2 = $targetPara.Characters($start2, 8)
$range2.Text = "effects "

Write-Host "Result:" $targetPara.Text
